# Weekly price-list update: a new daily record is inserted ahead of the
# existing row 59 (the whole table below it shifts down by one row,
# row 141 -> 142), and the new row is populated with its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 59; everything from the old row 59
# down to the old row 141 shifts down to rows 60..142.
$ws.Rows("59:59").Insert()

# Populate the newly inserted row 59 with the new record's data.
$ws.Range("A59").Value = 8
$ws.Range("B59").Value = "Terminal La Palmera de La Serena"
$ws.Range("C59").Value = "Coquimbo"
$ws.Range("D59").Value = 44679
$ws.Range("E59").Value = 4
$ws.Range("F59").Value = 100112040
$ws.Range("G59").Value = "Cilantro"
$ws.Range("H59").Value = "Sin especificar"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 2300
$ws.Range("K59").Value = 2500
$ws.Range("L59").Value = 3000
$ws.Range("M59").Value = 2750
$ws.Range("N59").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O59").Value = "Provincia del Elquí"
$ws.Range("P59").Value = 1833
$ws.Range("Q59").Value = 1.5
$ws.Range("R59").Value = "Hortaliza"
